# Commit: "Extract some repeated commands to single method"
# The accompanying test-data change adds a new expectation sheet,
# "verifyCheckoutOverview", holding the data that the extracted helper
# method is expected to verify on the cart-overview page, and nudges a
# couple of sheet selections left over from the author's last session.

$wb = $excel.ActiveWorkbook

# --- Leftover cursor/selection state on "verifyCartItemsCanBeRemoved" ---
$s7 = $wb.Worksheets.Item("verifyCartItemsCanBeRemoved")
$s7.Activate()
$s7.Range("A1:A4").Select() | Out-Null

# --- Leftover cursor/selection state on "verifyCheckoutInformation" ---
$s8 = $wb.Worksheets.Item("verifyCheckoutInformation")
$s8.Activate()
$s8.Range("G15").Select() | Out-Null

# --- Add the new "verifyCheckoutOverview" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$new.Name = "verifyCheckoutOverview"

$new.Range("A1").Value = "itemsToAdd"
$new.Range("A2").Value = "Sauce Labs Backpack"
$new.Range("A3").Value = "Sauce Labs Fleece Jacket, Sauce Labs Bolt T-Shirt, Test.allTheThings() T-Shirt (Red)"
$new.Range("A4").Value = "Sauce Labs Backpack, Sauce Labs Bike Light, Sauce Labs Bolt T-Shirt, Sauce Labs Onesie, Test.allTheThings() T-Shirt (Red)"

# Widen column A to fit the longest string, matching the sibling sheets.
$new.Columns.Item(1).ColumnWidth = 108.1666666667

# Leave the cursor where the author last left it on the new sheet.
$new.Range("C10").Select() | Out-Null
